$wb = $excel.ActiveWorkbook

# --- INTENT sheet: add a new intent "AMAZON.YesIntent" ---
$intent = $wb.Worksheets.Item("INTENT")
$intent.Cells.Item(10, 1).Value = "AMAZON.YesIntent"

# --- UTTERANCES_MAIN sheet: add a new column (F) for AMAZON.YesIntent ---
$utterances = $wb.Worksheets.Item("UTTERANCES_MAIN")
$utterances.Cells.Item(1, 6).Value = "AMAZON.YesIntent"
$utterances.Range("F2:F3").NumberFormat = "@"
$utterances.Cells.Item(2, 6).Value = "true"
$utterances.Cells.Item(3, 6).Value = "yes"

# --- Make UTTERANCES_MAIN the active sheet/tab ---
$utterances.Activate()
